# Daily update at 8 AM UTC
# Adds the next day's row of data (row 99) to the "Wins Over Time" sheet
# and moves the "last row" date-format style from the old last row (98)
# to the new last row (99).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98 was previously the last row and used the "last-row" date format
# (YYYY-MM-DD). Now that a new row follows it, it should use the regular
# date/time format used by all the other non-final rows.
$ws.Range("A98").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 99.
$ws.Range("A99").Value = 45838
$ws.Range("B99").Value = 420
$ws.Range("C99").Value = 416
$ws.Range("D99").Value = 430

# Row 99 is now the last row, so it gets the "last-row" date-only format.
$ws.Range("A99").NumberFormat = "YYYY-MM-DD"
